# Fix naming conventions for spr - add alt_label (column C) values
# for several rows in the variable-names sheet, and add a new shared
# string "f_spr" (note: only the shared string is added; the cell
# value used for row 268 points to it as per the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = alt_label. Fill in the alt_label for these rows to match
# their canonical/standardized short names.
$ws.Range("C203").Value = "spr"

$ws.Range("C260").Value = "fishing_mortality"
$ws.Range("C261").Value = "spr"
$ws.Range("C262").Value = "recruitment"
$ws.Range("C263").Value = "spawning_biomass"
$ws.Range("C264").Value = "biomass"
$ws.Range("C265").Value = "landings_weight"
$ws.Range("C266").Value = "discard_numbers"
$ws.Range("C267").Value = "discard_weight"
$ws.Range("C268").Value = "f_spr"

# Update the active selection to reflect where editing ended up.
$ws.Range("C269").Select()
